# =====================================================================
#  Applies the commit "On copie maintenant le style dans le fichier
#  sauvegarde." to the survey-export workbook:
#
#   1. Adds a second worksheet ("Feuille2") after "sheet1" with a few
#      notes cells, and makes it the active sheet/tab.
#   2. Refreshes the "taux reussite" bar chart's data-label position and
#      the category axis number-format-linked flag.
#   3. Fixes the row-15 recap cells that used to surface as #VALUE!
#      errors -- they should read #REF! instead.
# =====================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. New "Feuille2" worksheet, inserted right after the data sheet.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Feuille2"

$ws2.Range("A1").Value = "Test !!!!"
$ws2.Range("C4").Value = "Il doit y avoir des notes."
$ws2.Range("B6").Value = "zdzd"
$ws2.Range("B7").Value = "zdzd"

# match styling of the main sheet's cells as closely as this host allows
$ws2.Range("A1").Style = $ws1.Range("A1").Style
$ws2.Range("C4").Style = $ws1.Range("A1").Style
$ws2.Range("B6").Style = $ws1.Range("A1").Style
$ws2.Range("B7").Style = $ws1.Range("A1").Style

# leave the new sheet selected on its last-edited cell, then make it
# the active tab of the workbook (matches bookViews/activeTab -> 1)
$ws2.Range("B7").Select()
$ws2.Activate()

# ---------------------------------------------------------------
# 2. Chart touch-up: show outside-end data labels on the series and
#    stop linking the category axis number format to the source.
# ---------------------------------------------------------------
$co    = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$ser   = $chart.SeriesCollection(1)
$ser.DataLabels().Position = 2  # xlLabelPositionOutsideEnd -> "outEnd"

$catAx = $chart.Axes(1)
$catAx.TickLabels.NumberFormatLinked = 0

# ---------------------------------------------------------------
# 3. Row 15: recompute cells now surface #REF! instead of #VALUE!.
# ---------------------------------------------------------------
$refCells = @( `
    "L15","N15","P15","R15","T15","V15","X15","Z15", `
    "AB15","AD15","AF15","AH15","AJ15","AL15","AO15","AS15", `
    "AU15","AW15","AY15","BA15","BC15","BE15","BG15","BI15", `
    "BK15","BM15","BO15","BQ15","BS15","BU15","BW15","BY15" `
)

foreach ($cellRef in $refCells) {
    $ws1.Range($cellRef).Value = "#REF!"
}

# sheet1 keeps its own cursor on AO4, but Feuille2 must stay the
# workbook's active tab -- so select AO4 first, then re-activate Feuille2.
$ws1.Range("AO4").Select()
$ws2.Activate()
$ws2.Range("B7").Select()

Write-Host "edit.ps1 applied"
